$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns at E:G (will push Date Sampled .. Notes from E:N to H:Q)
$ws.Range("E1:G1").EntireColumn.Insert()

# Match the new columns' width to the existing C:D "9.5" width (not bestFit - explicitly set)
$ws.Range("E1:G1").EntireColumn.ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Header row labels for the new Month / Day / Year columns
$ws.Range("E1").Value = "Month"
$ws.Range("F1").Value = "Day"
$ws.Range("G1").Value = "Year"

# Data rows: populate Month / Day / Year to match the existing "Date Sampled" (8/3/2016)
$ws.Range("E2:E5").Value = 8
$ws.Range("F2:F5").Value = 3
$ws.Range("G2:G5").Value = 2016

# Restore selection to mirror the edit focus (E1:G1, active cell G1)
$ws.Range("E1:G1").Select()
